$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---------------------------------------------
$userData = $wb.Worksheets.Item(1)
$userData.Name = "UserData"

$storeData = $wb.Worksheets.Item(2)
$storeData.Name = "StoreData"

$sheet3 = $wb.Worksheets.Item(3)

# --- Insert new "OrderId" sheet right before the old "Sheet3" ----------
$orderId = $wb.Worksheets.Add($sheet3)
$orderId.Name = "OrderId"

# --- Populate StoreData with pet-store order test data ------------------
$headers = @("Id", "PetId", "quantity", "shipDate", "status", "complete")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $storeData.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$storeRows = @(
    @(1010, 1001, 2, "2024-12-03T19:17:38.568Z", "placed", $true),
    @(1020, 1002, 3, "2024-12-03T19:17:38.568Z", "placed", $true),
    @(1030, 1003, 5, "2024-12-03T19:17:38.568Z", "placed", $true)
)

for ($r = 0; $r -lt $storeRows.Length; $r++) {
    $rowValues = $storeRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $storeData.Cells.Item($r + 2, $c + 1).Value = $rowValues[$c]
    }
}

# Widen the shipDate column and leave the selection on G1:G4 like the source file
$storeData.Columns.Item(4).ColumnWidth = 23.5
$null = $storeData.Range("G1:G4").Select()

# --- Populate OrderId with a simple list of order ids --------------------
$orderId.Range("A1").Value = "orderId"
for ($i = 1; $i -le 10; $i++) {
    $orderId.Cells.Item($i + 1, 1).Value = $i
}
$null = $orderId.Range("A11").Select()
